$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 131.375
$ws.Range("I33").Value = 122.71429
$ws.Range("J33").Value = 192
$ws.Range("K33").Value = 122.71429
$ws.Range("L33").Value = 192
$ws.Range("M33").Value = 106.28571
$ws.Range("N33").Value = -650

$ws.Range("H53").Value = 273.5
$ws.Range("I53").Value = 87.75
$ws.Range("K53").Value = 87.75
$ws.Range("M53").Value = 549.25

$ws.Range("H100").Value = 1847.1875
$ws.Range("I100").Value = 1800.8334
$ws.Range("J100").Value = 1875
$ws.Range("K100").Value = 1800.8334
$ws.Range("L100").Value = 1875
$ws.Range("M100").Value = -1259.8334
$ws.Range("N100").Value = -2957

$ws.Range("H125").Value = 3368.0454
$ws.Range("I125").Value = 2185.0625
$ws.Range("J125").Value = 6522.6665
$ws.Range("K125").Value = 19665.5625
$ws.Range("L125").Value = 58703.9985
$ws.Range("M125").Value = -17205.5625
$ws.Range("N125").Value = -63623.9985

$ws.Range("H137").Value = 893.88
$ws.Range("I137").Value = 852.8421
$ws.Range("J137").Value = 1023.8333
$ws.Range("K137").Value = 2558.5263
$ws.Range("L137").Value = 3071.4999
$ws.Range("M137").Value = -8.526299999999992
$ws.Range("N137").Value = -8171.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 19800
$ws.Range("J23").Value = 19800
$ws.Range("L23").Value = 19800
$ws.Range("N23").Value = -20318

$ws.Range("H68").Value = 29000
$ws.Range("J68").Value = 29000
$ws.Range("L68").Value = 29000
$ws.Range("N68").Value = -30622

$ws.Range("H71").Value = 29000
$ws.Range("J71").Value = 29000
$ws.Range("L71").Value = 87000
$ws.Range("N71").Value = -95112

$ws.Range("H102").Value = 4829.8335
$ws.Range("I102").Value = 5326.3335
$ws.Range("K102").Value = 5326.3335
$ws.Range("M102").Value = -3704.3335

$ws.Range("H110").Value = 668.6667
$ws.Range("I110").Value = 668.6667
$ws.Range("K110").Value = 668.6667
$ws.Range("M110").Value = 1376.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 52633896
$ws.Range("I99").Value = 71430940
$ws.Range("J99").Value = 2164.4
$ws.Range("K99").Value = 71430940
$ws.Range("L99").Value = 2164.4
$ws.Range("M99").Value = -71429442
$ws.Range("N99").Value = -5160.4

$ws.Range("H105").Value = 8677.777
$ws.Range("I105").Value = 7728.5713
$ws.Range("K105").Value = 7728.5713
$ws.Range("M105").Value = -5981.5713

$ws.Range("H137").Value = 51667.617
$ws.Range("J137").Value = 51667.617
$ws.Range("L137").Value = 51667.617
$ws.Range("N137").Value = -61867.617

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 723.8333
$ws.Range("I16").Value = 606.7778
$ws.Range("J16").Value = 1075
$ws.Range("K16").Value = 606.7778
$ws.Range("L16").Value = 1075
$ws.Range("M16").Value = -319.7778
$ws.Range("N16").Value = -1649

$ws.Range("H31").Value = 3256.8
$ws.Range("I31").Value = 3408
$ws.Range("K31").Value = 3408
$ws.Range("M31").Value = -3113

$ws.Range("H34").Value = 3256.8
$ws.Range("I34").Value = 3408
$ws.Range("K34").Value = 3408
$ws.Range("M34").Value = -3206

$ws.Range("H113").Value = 723.8333
$ws.Range("I113").Value = 606.7778
$ws.Range("J113").Value = 1075
$ws.Range("K113").Value = 606.7778
$ws.Range("L113").Value = 1075
$ws.Range("M113").Value = 1563.2222
$ws.Range("N113").Value = -5415

$ws.Range("H140").Value = 55656
$ws.Range("J140").Value = 55656
$ws.Range("L140").Value = 55656
$ws.Range("N140").Value = -66016

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 421.66666
$ws.Range("I7").Value = 310
$ws.Range("J7").Value = 533.3333
$ws.Range("K7").Value = 930
$ws.Range("L7").Value = 1599.9999
$ws.Range("M7").Value = -818
$ws.Range("N7").Value = -1823.9999

$ws.Range("H80").Value = 5922.5557
$ws.Range("J80").Value = 5922.5557
$ws.Range("L80").Value = 17767.6671
$ws.Range("N80").Value = -19639.6671

$ws.Range("H83").Value = 5922.5557
$ws.Range("J83").Value = 5922.5557
$ws.Range("L83").Value = 53303.0013
$ws.Range("N83").Value = -62663.0013

$ws.Range("H87").Value = 10310
$ws.Range("J87").Value = 21575
$ws.Range("L87").Value = 64725
$ws.Range("N87").Value = -67221

$ws.Range("H90").Value = 10310
$ws.Range("J90").Value = 21575
$ws.Range("L90").Value = 194175
$ws.Range("N90").Value = -206655

$ws.Range("H92").Value = 320
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 8320
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 8320
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 8320
$ws.Range("N33").Value = -8824
$ws.Range("M33").ClearContents()

$ws.Range("H70").Value = 6370
$ws.Range("I70").Value = 5300
$ws.Range("K70").Value = 5300
$ws.Range("M70").Value = -5030

$ws.Range("H73").Value = 6370
$ws.Range("I73").Value = 5300
$ws.Range("K73").Value = 5300
$ws.Range("M73").Value = -4364

$ws.Range("H136").Value = 23192.334
$ws.Range("J136").Value = 23192.334
$ws.Range("L136").Value = 69577.00199999999
$ws.Range("N136").Value = -74677.00199999999

$ws.Range("H138").Value = 40900
$ws.Range("J138").Value = 40900
$ws.Range("L138").Value = 40900
$ws.Range("N138").Value = -51180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1915.8572
$ws.Range("I7").Value = 1818.5
$ws.Range("J7").Value = 2500
$ws.Range("K7").Value = 1818.5
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = -1706.5
$ws.Range("N7").Value = -2724

$ws.Range("H32").Value = 551.5
$ws.Range("I32").Value = 551.5
$ws.Range("K32").Value = 551.5
$ws.Range("M32").Value = -234.5

$ws.Range("H126").Value = 1915.8572
$ws.Range("I126").Value = 1818.5
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 5455.5
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2985.5
$ws.Range("N126").Value = -12440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 25000
$ws.Range("J39").Value = 25000
$ws.Range("L39").Value = 25000
$ws.Range("N39").Value = -25826

$ws.Range("H138").Value = 50760
$ws.Range("J138").Value = 50760
$ws.Range("L138").Value = 50760
$ws.Range("N138").Value = -61040
